$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Remove the review comment that was left on AA2 (author cleanup).
# ---------------------------------------------------------------------------
if ($ws.Range("AA2").Comment -ne $null) {
    $ws.Range("AA2").Comment.Delete()
}

# ---------------------------------------------------------------------------
# Row 2: update the VIN sample row from the old VOLKSWAGEN GOLF data to the
# new NISSAN ALTIMA data, and drop the BODY_STYLE_CD / BODYSHELL /
# ANTITHEFT_DISCOUNT / RESTRAINTS_DISCOUNT columns (K2, L2, AK2, AL2) that no
# longer apply to this record.
# ---------------------------------------------------------------------------
$ws.Range("K2").Clear()
$ws.Range("L2").Clear()
$ws.Range("AK2").Clear()
$ws.Range("AL2").Clear()

$ws.Range("A2").Value = "1N4BL3AP&H"
$ws.Range("B2").Value = "SYMBOL_2000_CHOICE"
$ws.Range("C2").Value = 2017
$ws.Range("D2").Value = "NISS"
$ws.Range("E2").Value = "NISSAN"
$ws.Range("F2").Value = "ALTIMA"
$ws.Range("G2").Value = "ALTIMA 3.5SL"
$ws.Range("H2").Value = 88888
$ws.Range("I2").Value = "SED"
$ws.Range("J2").Value = "SEDAN"
$ws.Range("M2").Value = "SED"
$ws.Range("N2").Value = "1.6L L4"
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = "G"
$ws.Range("Q2").Value = 214
$ws.Range("R2").Value = "2WD"
$ws.Range("S2").Value = 2
$ws.Range("T2").Value = 7
$ws.Range("U2").Value = "DUAL FRNT/HEAD/SIDE/REAR SIDE AIRBAGS"
$ws.Range("V2").Value = 2
$ws.Range("W2").Value = "4 WHEEL STANDARD"
$ws.Range("X2").Value = "STD"
$ws.Range("Y2").Value = "F-IMMOBILIZER/KEYLESS ENTRY"
$ws.Range("Z2").Value = "Q"
$ws.Range("AA2").Value = 20
$ws.Range("AB2").Value = 20
$ws.Range("AC2").Value = "B"
$ws.Range("AD2").Value = "N"
$ws.Range("AE2").Value = "B"
$ws.Range("AF2").Value = "B"
$ws.Range("AG2").Value = "B"
$ws.Range("AH2").Value = "B"
$ws.Range("AI2").Value = 20000101
$ws.Range("AJ2").Value = "N"

# ---------------------------------------------------------------------------
# Row 3: used to be a single blank, formatted cell (B3). Replace it with a
# second full VIN record (same shape as row 2) by first cloning row 2's
# formatting down, then filling in the new values.
# ---------------------------------------------------------------------------
$ws.Range("A2:AJ2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A3").Value = "1N4BL3AP&H"
$ws.Range("B3").Value = "SYMBOL_2000_CHOICE"
$ws.Range("C3").Value = 2017
$ws.Range("D3").Value = "NISS"
$ws.Range("E3").Value = "NISSAN MOTOR"
$ws.Range("F3").Value = "NISS ALTIMA"
$ws.Range("G3").Value = "ALTIMA 4.5SL"
$ws.Range("H3").Value = 88888
$ws.Range("I3").Value = "SED"
$ws.Range("J3").Value = "SED"
$ws.Range("K3").Clear()
$ws.Range("L3").Clear()
$ws.Range("M3").Value = "SED"
$ws.Range("N3").Value = "1.6L L4"
$ws.Range("O3").Value = 6
$ws.Range("P3").Value = "G"
$ws.Range("Q3").Value = 214
$ws.Range("R3").Value = "2WD"
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = 7
$ws.Range("U3").Value = "DUAL FRNT/HEAD/SIDE/REAR SIDE AIRBAGS"
$ws.Range("V3").Value = 2
$ws.Range("W3").Value = "4 WHEEL STANDARD"
$ws.Range("X3").Value = "STD"
$ws.Range("Y3").Value = "F-IMMOBILIZER/KEYLESS ENTRY"
$ws.Range("Z3").Value = "Q"
$ws.Range("AA3").Value = 20
$ws.Range("AB3").Value = 20
$ws.Range("AC3").Value = "B"
$ws.Range("AD3").Value = "N"
$ws.Range("AE3").Value = "B"
$ws.Range("AF3").Value = "B"
$ws.Range("AG3").Value = "B"
$ws.Range("AH3").Value = "B"
$ws.Range("AI3").Value = 20180612
$ws.Range("AJ3").Value = "Y"

# ---------------------------------------------------------------------------
# Move the active selection, mirroring the author re-selecting B9 before
# saving.
# ---------------------------------------------------------------------------
$ws.Range("B9").Select()
